# Applies the "Answered question 3 and 7, added a bit of hardware
# requirements, and use cases." edit to the Systems Requirements Document.

$d = $word.ActiveDocument

function Set-StdFont($rng) {
    # Re-apply the document's standard body formatting (Times New Roman,
    # 12pt / sz 24) to a range so newly-created runs match the rest of
    # the document.
    $rng.Font.Name = "Times New Roman"
    $rng.Font.NameAscii = "Times New Roman"
    $rng.Font.NameOther = "Times New Roman"
    $rng.Font.NameBi = "Times New Roman"
    $rng.Font.Size = 12
}

# ---------------------------------------------------------------------
# 1) Question 2 answer: flesh out the clerk inventory answer.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "2: Clerks will need to be able to view the inventory and be able to see how many of what game are in stock, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2: Clerks will need to be able to view and edit the inventory of the local store. Sales and purchases should also automatically update the store inventory.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) Question 3 answer.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "3: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3: Customers should be able to view store inventory on a screen local to the store but they should not be able to do anything other than view.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3) Question 7 answer.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "7: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "7: Managers need the highest clearance, Clerks need general access, and customers should have a view only screen.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 4) Use Case section: Manager / Clerk / Customer descriptions.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Manager: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Manager: Order games, manage stock, manage prices",
    2) | Out-Null

$d.Content.Find.Execute(
    "Clerk: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Clerk: Add and remove games from the inventory, view store inventory and game information.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Customer: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Customer: View store inventory.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 5) Hardware Requirements: fill the blank paragraph and add two more.
# ---------------------------------------------------------------------
$hwHeading = $d.Content.Find
$hwPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Hardware Requirements`r") {
        $hwPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

$hwRange = $hwPara.Range
$hwRange.Text = "Customer will need access to a local touch screen."
$hwRange2 = $hwPara.Range
Set-StdFont $hwRange2

$hwEnd = $hwPara.Range
$insertPoint1 = $d.Range($hwEnd.End - 1, $hwEnd.End - 1)
$insertPoint1.InsertParagraphAfter()
$clerkHwPara = $d.Paragraphs.Item($hwPara.Range.Information(3) + 1)

# Re-locate paragraphs by re-scanning, since indices can shift.
$hwIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Customer will need access to a local touch screen.`r") {
        $hwIndex = $i
        break
    }
}

$newPara1 = $d.Paragraphs.Item($hwIndex + 1)
$newPara1.Range.Text = "Clerks will need access to a local PC including monitor and keyboard"
Set-StdFont ($d.Paragraphs.Item($hwIndex + 1).Range)

$insertPoint2Para = $d.Paragraphs.Item($hwIndex + 1)
$ip2 = $d.Range($insertPoint2Para.Range.End - 1, $insertPoint2Para.Range.End - 1)
$ip2.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($hwIndex + 2)
$newPara2.Range.Text = "Managers will have access to the central database"
Set-StdFont ($d.Paragraphs.Item($hwIndex + 2).Range)

# ---------------------------------------------------------------------
# 6) Non-Functional Requirements: add a closing sentence.
# ---------------------------------------------------------------------
$nfIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "*Knowledge of the user?`r") {
        $nfIndex = $i
        break
    }
}
$nfPara = $d.Paragraphs.Item($nfIndex)
$ip3 = $d.Range($nfPara.Range.End - 1, $nfPara.Range.End - 1)
$ip3.InsertParagraphAfter()
$newPara3 = $d.Paragraphs.Item($nfIndex + 1)
$newPara3.Range.Text = "Clerks and Managers will need to be able to use windows and have a basic understanding of how to use a PC. "
Set-StdFont ($d.Paragraphs.Item($nfIndex + 1).Range)

Write-Host "Done"
